$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48
$ws.Cells.Item($row, 1).Value = "Team Pesto nel Tigullio"
$ws.Cells.Item($row, 2).Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Cells.Item($row, 3).Value = "Federico  Manica | iMontagna"
$ws.Cells.Item($row, 4).Value = "Alessandro Comper | F.C. Gorillaz"
$ws.Cells.Item($row, 5).Value = "Michele Merighi | Clitoriders"
$ws.Cells.Item($row, 6).Value = "Maverick  Bertolini | A.C. Denti"
